# "Loocal phone for tests" - swap the test recipient's phone number for a
# local (Guinea, +224) number, and touch up the header row's look so it
# matches the other locale's template (explicit black text instead of the
# theme color, and a slightly shorter header row height).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Core edit: replace the phone number shared string used by C2.
$ws.Range("C2").Value = "+224 628 11 11 86"

# Cosmetic touch-up that came along with the resave: the bordered header
# cells (B1, C1, A2) get an explicit black font color instead of the
# automatic/theme color.
$ws.Range("B1").Font.Color = 0
$ws.Range("C1").Font.Color = 0
$ws.Range("A2").Font.Color = 0

# Header rows shrink slightly (20.25pt -> 19.5pt).
$ws.Rows.Item(1).RowHeight = 19.5
$ws.Rows.Item(2).RowHeight = 19.5
